# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.311.18'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.017.62'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.82'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.98'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.015.84'
$ws.Range('E8').Value = '  +0.72%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  +10.88%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.51'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').Value = '3.516.84'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.02'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '62.250.13'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '3.021.02'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '448.75'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.42'
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.34'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.90'
$ws.Range('E25').Value = '  +12.03%  '
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.05'
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.71'
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  +3.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.14'
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.55'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('E35').Value = '  +6.10%  '
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +5.23%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.15'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.286'
$ws.Range('E43').Value = '  +7.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.06'
$ws.Range('E44').Value = '  +9.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '394.91'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '2.738.22'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.85'
$ws.Range('E48').Value = '  +4.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.19'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('E51').Value = '  -1.02%  '
